# Applies the cryptos-list data refresh described by the commit diff.
# D-column values are numeric-looking text (e.g. "1.001", "29.913.71") that
# must stay as literal text, matching the source inlineStr cells -- so the
# touched D cells are pre-formatted as Text before the values are written,
# exactly as Excel requires to avoid silently re-typing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The whole Price column already holds plain-text values (e.g. "29.907.47",
# "1.002"); keep that convention by formatting it as Text before writing the
# refreshed figures, so Excel doesn't silently re-type them as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.913.71"
$ws.Range("E2").Value = "  -1.17%  "

# Row 3
$ws.Range("D3").Value = "1.918.87"
$ws.Range("E3").Value = "  +1.31%  "

# Row 4
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "320.32"
$ws.Range("E5").Value = "  -1.23%  "

# Row 6
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("E7").Value = "  -2.44%  "

# Row 8
$ws.Range("D8").Value = "0.4029"
$ws.Range("E8").Value = "  +0.26%  "

# Row 9
$ws.Range("D9").Value = "0.08322"
$ws.Range("E9").Value = "  -1.07%  "

# Row 10
$ws.Range("D10").Value = "42.70"
$ws.Range("E10").Value = "  +0.01%  "

# Row 11
$ws.Range("D11").Value = "1.103"
$ws.Range("E11").Value = "  -1.10%  "

# Row 12
$ws.Range("D12").Value = "23.75"
$ws.Range("E12").Value = "  +2.48%  "

# Row 13
$ws.Range("D13").Value = "1.920.95"
$ws.Range("E13").Value = "  +1.27%  "

# Row 14
$ws.Range("D14").Value = "6.401"
$ws.Range("E14").Value = "  -0.60%  "

# Row 15
$ws.Range("D15").Value = "7.223"
$ws.Range("E15").Value = "  -1.30%  "

# Row 16
$ws.Range("D16").Value = "1.000"
$ws.Range("E16").Value = "  -0.12%  "

# Row 17
$ws.Range("D17").Value = "92.13"
$ws.Range("E17").Value = "  -2.33%  "

# Row 18
$ws.Range("E18").Value = "  -1.14%  "

# Row 19
$ws.Range("D19").Value = "0.06508"
$ws.Range("E19").Value = "  -2.03%  "

# Row 20
$ws.Range("D20").Value = "18.26"
$ws.Range("E20").Value = "  +0.15%  "

# Row 21
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22
$ws.Range("D22").Value = "5.943"
$ws.Range("E22").Value = "  -0.23%  "

# Row 23
$ws.Range("D23").Value = "29.944.64"
$ws.Range("E23").Value = "  -1.01%  "

# Row 24
$ws.Range("D24").Value = "11.33"
$ws.Range("E24").Value = "  +0.31%  "

# Row 25
$ws.Range("D25").Value = "2.192"
$ws.Range("E25").Value = "  -1.66%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "22.15"
$ws.Range("E26").Value = "  +2.23%  "

# Row 27
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.142.34"
$ws.Range("E27").Value = "  +1.31%  "

# Row 28
$ws.Range("D28").Value = "162.01"
$ws.Range("E28").Value = "  +0.10%  "

# Row 29
$ws.Range("D29").Value = "2.317"
$ws.Range("E29").Value = "  -1.23%  "

# Row 30
$ws.Range("D30").Value = "128.96"
$ws.Range("E30").Value = "  -0.26%  "

# Row 31
$ws.Range("D31").Value = "1.130"
$ws.Range("E31").Value = "  +3.66%  "

# Row 32
$ws.Range("E32").Value = "  -1.76%  "

# Row 33
$ws.Range("D33").Value = "5.959"
$ws.Range("E33").Value = "  -2.26%  "

# Row 34
$ws.Range("D34").Value = "3.792"
$ws.Range("E34").Value = "  +1.34%  "

# Row 35
$ws.Range("E35").Value = "  -1.78%  "

# Row 36
$ws.Range("D36").Value = "5.395"
$ws.Range("E36").Value = "  +1.14%  "

# Row 37
$ws.Range("D37").Value = "0.06411"
$ws.Range("E37").Value = "  -2.17%  "

# Row 38
$ws.Range("D38").Value = "0.2155"
$ws.Range("E38").Value = "  -1.94%  "

# Row 39
$ws.Range("D39").Value = "0.6502"
$ws.Range("E39").Value = "  +0.04%  "

# Row 40
$ws.Range("D40").Value = "8.723"
$ws.Range("E40").Value = "  -1.11%  "

# Row 41
$ws.Range("E41").Value = "  -2.21%  "

# Row 42
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "1.218"
$ws.Range("E42").Value = "  -0.90%  "

# Row 43
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "11.38"
$ws.Range("E43").Value = "  -3.63%  "

# Row 44
$ws.Range("D44").Value = "2.235"
$ws.Range("E44").Value = "  +8.77%  "

# Row 45
$ws.Range("D45").Value = "13.48"
$ws.Range("E45").Value = "  +1.51%  "

# Row 46
$ws.Range("D46").Value = "0.6096"
$ws.Range("E46").Value = "  +0.19%  "

# Row 47
$ws.Range("D47").Value = "3.638"
$ws.Range("E47").Value = "  -1.25%  "

# Row 48
$ws.Range("D48").Value = "1.209"
$ws.Range("E48").Value = "  -2.21%  "

# Row 49
$ws.Range("D49").Value = "121.90"
$ws.Range("E49").Value = "  -2.21%  "

# Row 50
$ws.Range("D50").Value = "78.94"
$ws.Range("E50").Value = "  -0.22%  "

# Row 51
$ws.Range("D51").Value = "1.128"

